# Trade #21 (MarketMaking strategy) closed at 2026-02-17 20:53:30 - unknown UNKNOWN +0.000%
# and a brand new trade (#82) opened at 2026-02-17 20:53:24.
# This updates the Summary, Strategy Status, All Trades and MarketMaking
# sheets to reflect both events.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - headline stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 49      # Total Trades
$summary.Range("B9").Value = 44.9    # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 16       # Trades
$status.Range("G5").Value = 50       # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet - close out trade #49 (row 50) and append trade #82
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close trade #49
$allTrades.Range("G50").Value = 0.82
$allTrades.Range("H50").Value = "CLOSED"
$allTrades.Range("K50").Value = 100.4
$allTrades.Range("L50").Value = "early_exit"
$allTrades.Range("M50").Value = 0.13

# Append new trade #82 as row 83
$allTrades.Range("A83").Value = 82
$allTrades.Range("B83").Value = "'2026-02-17"
$allTrades.Range("B83").Style = "Normal"
$allTrades.Range("C83").Value = "'20:53:24"
$allTrades.Range("C83").Style = "Normal"
$allTrades.Range("D83").Value = "MarketMaking"
$allTrades.Range("E83").Value = "DOWN"
$allTrades.Range("F83").Value = 0.82
$allTrades.Range("H83").Value = "OPEN"
$allTrades.Range("I83").Value = 0
$allTrades.Range("J83").Value = 0
$allTrades.Range("K83").Value = 100.3984370824165
$allTrades.Range("M83").Value = 0
$allTrades.Range("N83").Value = 0
$allTrades.Range("O83").Value = 0
$allTrades.Range("P83").Value = 0.6
$allTrades.Range("Q83").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet - close out trade #49 (row 17) and append trade #82
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close trade #49
$mm.Range("G17").Value = 0.82
$mm.Range("H17").Value = "CLOSED"
$mm.Range("K17").Value = 100.4
$mm.Range("P17").Value = "early_exit"
$mm.Range("Q17").Value = 0.13

# Append new trade #82 as row 50
$mm.Range("A50").Value = 82
$mm.Range("B50").Value = "'2026-02-17"
$mm.Range("B50").Style = "Normal"
$mm.Range("C50").Value = "'20:53:24"
$mm.Range("C50").Style = "Normal"
$mm.Range("D50").Value = "MarketMaking"
$mm.Range("E50").Value = "DOWN"
$mm.Range("F50").Value = 0.82
$mm.Range("H50").Value = "OPEN"
$mm.Range("I50").Value = 0
$mm.Range("J50").Value = 0
$mm.Range("K50").Value = 100.3984370824165
$mm.Range("L50").Value = 0
$mm.Range("M50").Value = 0
$mm.Range("N50").Value = 0.6
$mm.Range("O50").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q50").Value = 0
